$wb = $excel.ActiveWorkbook

$debts = $wb.Worksheets.Add()
$debts.Name = "Debts"

$fixed = $wb.Worksheets.Add()
$fixed.Name = "Fixed Assets"
